$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value2 = 'negative'
$ws.Range('J1').Value2 = 'positive'

$ws.Range('A2').Value2 = 'name'
$ws.Range('B2').Value2 = 'anchor score'
$ws.Range('C2').Value2 = 'type occurences'
$ws.Range('D2').Value2 = 'total occurences'
$ws.Range('E2').Value2 = '+%'
$ws.Range('F2').Value2 = '-%'
$ws.Range('G2').Value2 = 'both'
$ws.Range('H2').Value2 = 'normal'
$ws.Range('J2').Value2 = 'name'
$ws.Range('K2').Value2 = 'anchor score'
$ws.Range('L2').Value2 = 'type occurences'
$ws.Range('M2').Value2 = 'total occurences'
$ws.Range('N2').Value2 = '+%'
$ws.Range('O2').Value2 = '-%'
$ws.Range('P2').Value2 = 'both'
$ws.Range('Q2').Value2 = 'normal'

$ws.Range('A3').Value2 = 'poorly'
$ws.Range('B3').Value2 = 0.9565217391304348
$ws.Range('C3').Value2 = 44
$ws.Range('D3').Value2 = 44
$ws.Range('E3').Value2 = 0
$ws.Range('F3').Value2 = 1
$ws.Range('G3').Value2 = $False
$ws.Range('H3').Value2 = 2
$ws.Range('J3').Value2 = 'wonderful'
$ws.Range('K3').Value2 = 0.8928571428571429
$ws.Range('L3').Value2 = 50
$ws.Range('M3').Value2 = 50
$ws.Range('N3').Value2 = 1
$ws.Range('O3').Value2 = 0
$ws.Range('P3').Value2 = $False
$ws.Range('Q3').Value2 = 6

$ws.Range('A4').Value2 = 'disappointing'
$ws.Range('B4').Value2 = 0.8863636363636364
$ws.Range('C4').Value2 = 39
$ws.Range('D4').Value2 = 39
$ws.Range('E4').Value2 = 0
$ws.Range('F4').Value2 = 1
$ws.Range('G4').Value2 = $False
$ws.Range('H4').Value2 = 5
$ws.Range('J4').Value2 = 'awesome'
$ws.Range('K4').Value2 = 0.8923076923076924
$ws.Range('L4').Value2 = 58
$ws.Range('M4').Value2 = 58
$ws.Range('N4').Value2 = 1
$ws.Range('O4').Value2 = 0
$ws.Range('P4').Value2 = $False
$ws.Range('Q4').Value2 = 7

$ws.Range('A5').Value2 = 'thin'
$ws.Range('B5').Value2 = 0.7586206896551724
$ws.Range('C5').Value2 = 22
$ws.Range('D5').Value2 = 22
$ws.Range('E5').Value2 = 0
$ws.Range('F5').Value2 = 1
$ws.Range('G5').Value2 = $False
$ws.Range('H5').Value2 = 7
$ws.Range('J5').Value2 = 'amazing'
$ws.Range('K5').Value2 = 0.8888888888888888
$ws.Range('L5').Value2 = 24
$ws.Range('M5').Value2 = 24
$ws.Range('N5').Value2 = 1
$ws.Range('O5').Value2 = 0
$ws.Range('P5').Value2 = $False
$ws.Range('Q5').Value2 = 3

$ws.Range('A6').Value2 = 'returned'
$ws.Range('B6').Value2 = 0.7368421052631579
$ws.Range('C6').Value2 = 28
$ws.Range('D6').Value2 = 28
$ws.Range('E6').Value2 = 0
$ws.Range('F6').Value2 = 1
$ws.Range('G6').Value2 = $False
$ws.Range('H6').Value2 = 10
$ws.Range('J6').Value2 = 'favorite'
$ws.Range('K6').Value2 = 0.8494623655913979
$ws.Range('L6').Value2 = 79
$ws.Range('M6').Value2 = 79
$ws.Range('N6').Value2 = 1
$ws.Range('O6').Value2 = 0
$ws.Range('P6').Value2 = $False
$ws.Range('Q6').Value2 = 14

$ws.Range('A7').Value2 = 'poor'
$ws.Range('B7').Value2 = 0.704225352112676
$ws.Range('C7').Value2 = 50
$ws.Range('D7').Value2 = 50
$ws.Range('E7').Value2 = 0
$ws.Range('F7').Value2 = 1
$ws.Range('G7').Value2 = $False
$ws.Range('H7').Value2 = 21
$ws.Range('J7').Value2 = 'excellent'
$ws.Range('K7').Value2 = 0.71875
$ws.Range('L7').Value2 = 46
$ws.Range('M7').Value2 = 46
$ws.Range('N7').Value2 = 1
$ws.Range('O7').Value2 = 0
$ws.Range('P7').Value2 = $False
$ws.Range('Q7').Value2 = 18

$ws.Range('A8').Value2 = 'however'
$ws.Range('B8').Value2 = 0.703125
$ws.Range('C8').Value2 = 45
$ws.Range('D8').Value2 = 45
$ws.Range('E8').Value2 = 0
$ws.Range('F8').Value2 = 1
$ws.Range('G8').Value2 = $False
$ws.Range('H8').Value2 = 19
$ws.Range('J8').Value2 = 'classic'
$ws.Range('K8').Value2 = 0.660377358490566
$ws.Range('L8').Value2 = 35
$ws.Range('M8').Value2 = 35
$ws.Range('N8').Value2 = 1
$ws.Range('O8').Value2 = 0
$ws.Range('P8').Value2 = $False
$ws.Range('Q8').Value2 = 18

$ws.Range('A9').Value2 = 'disappointed'
$ws.Range('B9').Value2 = 0.6881720430107527
$ws.Range('C9').Value2 = 128
$ws.Range('D9').Value2 = 128
$ws.Range('E9').Value2 = 0
$ws.Range('F9').Value2 = 1
$ws.Range('G9').Value2 = $False
$ws.Range('H9').Value2 = 58
$ws.Range('J9').Value2 = 'love'
$ws.Range('K9').Value2 = 0.5796269727403156
$ws.Range('L9').Value2 = 404
$ws.Range('M9').Value2 = 404
$ws.Range('N9').Value2 = 1
$ws.Range('O9').Value2 = 0
$ws.Range('P9').Value2 = $False
$ws.Range('Q9').Value2 = 293

$ws.Range('A10').Value2 = 'junk'
$ws.Range('B10').Value2 = 0.6727272727272727
$ws.Range('C10').Value2 = 37
$ws.Range('D10').Value2 = 37
$ws.Range('E10').Value2 = 0
$ws.Range('F10').Value2 = 1
$ws.Range('G10').Value2 = $False
$ws.Range('H10').Value2 = 18
$ws.Range('J10').Value2 = 'thank'
$ws.Range('K10').Value2 = 0.5362318840579711
$ws.Range('L10').Value2 = 37
$ws.Range('M10').Value2 = 37
$ws.Range('N10').Value2 = 1
$ws.Range('O10').Value2 = 0
$ws.Range('P10').Value2 = $False
$ws.Range('Q10').Value2 = 32

$ws.Range('A11').Value2 = 'broke'
$ws.Range('B11').Value2 = 0.6699029126213593
$ws.Range('C11').Value2 = 138
$ws.Range('D11').Value2 = 138
$ws.Range('E11').Value2 = 0
$ws.Range('F11').Value2 = 1
$ws.Range('G11').Value2 = $False
$ws.Range('H11').Value2 = 68
$ws.Range('J11').Value2 = 'loves'
$ws.Range('K11').Value2 = 0.504149377593361
$ws.Range('L11').Value2 = 243
$ws.Range('M11').Value2 = 243
$ws.Range('N11').Value2 = 1
$ws.Range('O11').Value2 = 0
$ws.Range('P11').Value2 = $False
$ws.Range('Q11').Value2 = 239

$ws.Range('A12').Value2 = 'waste'
$ws.Range('B12').Value2 = 0.6351351351351351
$ws.Range('C12').Value2 = 94
$ws.Range('D12').Value2 = 94
$ws.Range('E12').Value2 = 0
$ws.Range('F12').Value2 = 1
$ws.Range('G12').Value2 = $False
$ws.Range('H12').Value2 = 54
$ws.Range('J12').Value2 = 'great'
$ws.Range('K12').Value2 = 0.4721311475409836
$ws.Range('L12').Value2 = 576
$ws.Range('M12').Value2 = 576
$ws.Range('N12').Value2 = 1
$ws.Range('O12').Value2 = 0
$ws.Range('P12').Value2 = $False
$ws.Range('Q12').Value2 = 644

$ws.Range('A13').Value2 = 'tiny'
$ws.Range('B13').Value2 = 0.6285714285714286
$ws.Range('C13').Value2 = 22
$ws.Range('D13').Value2 = 22
$ws.Range('E13').Value2 = 0
$ws.Range('F13').Value2 = 1
$ws.Range('G13').Value2 = $False
$ws.Range('H13').Value2 = 13
$ws.Range('J13').Value2 = 'pleased'
$ws.Range('K13').Value2 = 0.4150943396226415
$ws.Range('L13').Value2 = 22
$ws.Range('M13').Value2 = 22
$ws.Range('N13').Value2 = 1
$ws.Range('O13').Value2 = 0
$ws.Range('P13').Value2 = $False
$ws.Range('Q13').Value2 = 31

$ws.Range('A14').Value2 = 'water'
$ws.Range('B14').Value2 = 0.6190476190476191
$ws.Range('C14').Value2 = 26
$ws.Range('D14').Value2 = 26
$ws.Range('E14').Value2 = 0
$ws.Range('F14').Value2 = 1
$ws.Range('G14').Value2 = $False
$ws.Range('H14').Value2 = 16
$ws.Range('J14').Value2 = 'loved'
$ws.Range('K14').Value2 = 0.382262996941896
$ws.Range('L14').Value2 = 125
$ws.Range('M14').Value2 = 125
$ws.Range('N14').Value2 = 1
$ws.Range('O14').Value2 = 0
$ws.Range('P14').Value2 = $False
$ws.Range('Q14').Value2 = 202

$ws.Range('A15').Value2 = 'smaller'
$ws.Range('B15').Value2 = 0.6050420168067226
$ws.Range('C15').Value2 = 72
$ws.Range('D15').Value2 = 72
$ws.Range('E15').Value2 = 0
$ws.Range('F15').Value2 = 1
$ws.Range('G15').Value2 = $False
$ws.Range('H15').Value2 = 47
$ws.Range('J15').Value2 = 'perfect'
$ws.Range('K15').Value2 = 0.3734939759036144
$ws.Range('L15').Value2 = 62
$ws.Range('M15').Value2 = 62
$ws.Range('N15').Value2 = 1
$ws.Range('O15').Value2 = 0
$ws.Range('P15').Value2 = $False
$ws.Range('Q15').Value2 = 104

$ws.Range('A16').Value2 = 'instead'
$ws.Range('B16').Value2 = 0.5833333333333334
$ws.Range('C16').Value2 = 28
$ws.Range('D16').Value2 = 28
$ws.Range('E16').Value2 = 0
$ws.Range('F16').Value2 = 1
$ws.Range('G16').Value2 = $False
$ws.Range('H16').Value2 = 20
$ws.Range('J16').Value2 = 'friends'
$ws.Range('K16').Value2 = 0.3333333333333333
$ws.Range('L16').Value2 = 63
$ws.Range('M16').Value2 = 63
$ws.Range('N16').Value2 = 1
$ws.Range('O16').Value2 = 0
$ws.Range('P16').Value2 = $False
$ws.Range('Q16').Value2 = 126

$ws.Range('A17').Value2 = 'small'
$ws.Range('B17').Value2 = 0.5246376811594203
$ws.Range('C17').Value2 = 181
$ws.Range('D17').Value2 = 181
$ws.Range('E17').Value2 = 0
$ws.Range('F17').Value2 = 1
$ws.Range('G17').Value2 = $False
$ws.Range('H17').Value2 = 164
$ws.Range('J17').Value2 = 'best'
$ws.Range('K17').Value2 = 0.325
$ws.Range('L17').Value2 = 39
$ws.Range('M17').Value2 = 39
$ws.Range('N17').Value2 = 1
$ws.Range('O17').Value2 = 0
$ws.Range('P17').Value2 = $False
$ws.Range('Q17').Value2 = 81

$ws.Range('A18').Value2 = 'guess'
$ws.Range('B18').Value2 = 0.5
$ws.Range('C18').Value2 = 27
$ws.Range('D18').Value2 = 27
$ws.Range('E18').Value2 = 0
$ws.Range('F18').Value2 = 1
$ws.Range('G18').Value2 = $False
$ws.Range('H18').Value2 = 27
$ws.Range('J18').Value2 = 'learn'
$ws.Range('K18').Value2 = 0.234375
$ws.Range('L18').Value2 = 30
$ws.Range('M18').Value2 = 30
$ws.Range('N18').Value2 = 1
$ws.Range('O18').Value2 = 0
$ws.Range('P18').Value2 = $False
$ws.Range('Q18').Value2 = 98

$ws.Range('A19').Value2 = 'broken'
$ws.Range('B19').Value2 = 0.4698795180722892
$ws.Range('C19').Value2 = 39
$ws.Range('D19').Value2 = 39
$ws.Range('E19').Value2 = 0
$ws.Range('F19').Value2 = 1
$ws.Range('G19').Value2 = $False
$ws.Range('H19').Value2 = 44
$ws.Range('J19').Value2 = 'enjoyed'
$ws.Range('K19').Value2 = 0.2301587301587301
$ws.Range('L19').Value2 = 29
$ws.Range('M19').Value2 = 29
$ws.Range('N19').Value2 = 1
$ws.Range('O19').Value2 = 0
$ws.Range('P19').Value2 = $False
$ws.Range('Q19').Value2 = 97

$ws.Range('A20').Value2 = 'plastic'
$ws.Range('B20').Value2 = 0.4566929133858268
$ws.Range('C20').Value2 = 58
$ws.Range('D20').Value2 = 58
$ws.Range('E20').Value2 = 0
$ws.Range('F20').Value2 = 1
$ws.Range('G20').Value2 = $False
$ws.Range('H20').Value2 = 69
$ws.Range('J20').Value2 = 'christmas'
$ws.Range('K20').Value2 = 0.2168674698795181
$ws.Range('L20').Value2 = 54
$ws.Range('M20').Value2 = 54
$ws.Range('N20').Value2 = 1
$ws.Range('O20').Value2 = 0
$ws.Range('P20').Value2 = $False
$ws.Range('Q20').Value2 = 195

$ws.Range('A21').Value2 = 'pay'
$ws.Range('B21').Value2 = 0.4285714285714285
$ws.Range('C21').Value2 = 27
$ws.Range('D21').Value2 = 27
$ws.Range('E21').Value2 = 0
$ws.Range('F21').Value2 = 1
$ws.Range('G21').Value2 = $False
$ws.Range('H21').Value2 = 36
$ws.Range('J21').Value2 = 'enjoy'
$ws.Range('K21').Value2 = 0.2150537634408602
$ws.Range('L21').Value2 = 40
$ws.Range('M21').Value2 = 40
$ws.Range('N21').Value2 = 1
$ws.Range('O21').Value2 = 0
$ws.Range('P21').Value2 = $False
$ws.Range('Q21').Value2 = 146

$ws.Range('A22').Value2 = 'apart'
$ws.Range('B22').Value2 = 0.4105263157894737
$ws.Range('C22').Value2 = 39
$ws.Range('D22').Value2 = 39
$ws.Range('E22').Value2 = 0
$ws.Range('F22').Value2 = 1
$ws.Range('G22').Value2 = $False
$ws.Range('H22').Value2 = 56
$ws.Range('J22').Value2 = 'happy'
$ws.Range('K22').Value2 = 0.2097902097902098
$ws.Range('L22').Value2 = 30
$ws.Range('M22').Value2 = 30
$ws.Range('N22').Value2 = 1
$ws.Range('O22').Value2 = 0
$ws.Range('P22').Value2 = $False
$ws.Range('Q22').Value2 = 113

$ws.Range('A23').Value2 = 'ok'
$ws.Range('B23').Value2 = 0.3984375
$ws.Range('C23').Value2 = 51
$ws.Range('D23').Value2 = 51
$ws.Range('E23').Value2 = 0
$ws.Range('F23').Value2 = 1
$ws.Range('G23').Value2 = $False
$ws.Range('H23').Value2 = 77
$ws.Range('J23').Value2 = 'fun'
$ws.Range('K23').Value2 = 0.1747146619841967
$ws.Range('L23').Value2 = 199
$ws.Range('M23').Value2 = 201
$ws.Range('N23').Value2 = 0.99
$ws.Range('O23').Value2 = 0.01000000000000001
$ws.Range('P23').Value2 = $True
$ws.Range('Q23').Value2 = 940

$ws.Range('A24').Value2 = 'di'
$ws.Range('B24').Value2 = 0.375
$ws.Range('C24').Value2 = 24
$ws.Range('D24').Value2 = 24
$ws.Range('E24').Value2 = 0
$ws.Range('F24').Value2 = 1
$ws.Range('G24').Value2 = $False
$ws.Range('H24').Value2 = 40
$ws.Range('J24').Value2 = 'family'
$ws.Range('K24').Value2 = 0.1002785515320334
$ws.Range('L24').Value2 = 36
$ws.Range('M24').Value2 = 36
$ws.Range('N24').Value2 = 1
$ws.Range('O24').Value2 = 0
$ws.Range('P24').Value2 = $False
$ws.Range('Q24').Value2 = 323

$ws.Range('A25').Value2 = 'cheap'
$ws.Range('B25').Value2 = 0.3744075829383886
$ws.Range('C25').Value2 = 79
$ws.Range('D25').Value2 = 79
$ws.Range('E25').Value2 = 0
$ws.Range('F25').Value2 = 1
$ws.Range('G25').Value2 = $False
$ws.Range('H25').Value2 = 132
$ws.Range('J25').Value2 = 'easy'
$ws.Range('K25').Value2 = 0.09946236559139784
$ws.Range('L25').Value2 = 37
$ws.Range('M25').Value2 = 39
$ws.Range('N25').Value2 = 0.95
$ws.Range('O25').Value2 = 0.05000000000000004
$ws.Range('P25').Value2 = $True
$ws.Range('Q25').Value2 = 335

$ws.Range('A26').Value2 = 'paint'
$ws.Range('B26').Value2 = 0.3650793650793651
$ws.Range('C26').Value2 = 23
$ws.Range('D26').Value2 = 23
$ws.Range('E26').Value2 = 0
$ws.Range('F26').Value2 = 1
$ws.Range('G26').Value2 = $False
$ws.Range('H26').Value2 = 40
$ws.Range('J26').Value2 = 'game'
$ws.Range('K26').Value2 = 0.09935064935064936
$ws.Range('L26').Value2 = 153
$ws.Range('M26').Value2 = 154
$ws.Range('N26').Value2 = 0.99
$ws.Range('O26').Value2 = 0.01000000000000001
$ws.Range('P26').Value2 = $True
$ws.Range('Q26').Value2 = 1387

$ws.Range('A27').Value2 = 'thought'
$ws.Range('B27').Value2 = 0.2722772277227723
$ws.Range('C27').Value2 = 55
$ws.Range('D27').Value2 = 55
$ws.Range('E27').Value2 = 0
$ws.Range('F27').Value2 = 1
$ws.Range('G27').Value2 = $False
$ws.Range('H27').Value2 = 147
$ws.Range('J27').Value2 = 'play'
$ws.Range('K27').Value2 = 0.04933333333333333
$ws.Range('L27').Value2 = 37
$ws.Range('M27').Value2 = 39
$ws.Range('N27').Value2 = 0.95
$ws.Range('O27').Value2 = 0.05000000000000004
$ws.Range('P27').Value2 = $True
$ws.Range('Q27').Value2 = 713

$ws.Range('A28').Value2 = 'bit'
$ws.Range('B28').Value2 = 0.2653061224489796
$ws.Range('C28').Value2 = 26
$ws.Range('D28').Value2 = 26
$ws.Range('E28').Value2 = 0
$ws.Range('F28').Value2 = 1
$ws.Range('G28').Value2 = $False
$ws.Range('H28').Value2 = 72

$ws.Range('A29').Value2 = 'though'
$ws.Range('B29').Value2 = 0.264957264957265
$ws.Range('C29').Value2 = 31
$ws.Range('D29').Value2 = 31
$ws.Range('E29').Value2 = 0
$ws.Range('F29').Value2 = 1
$ws.Range('G29').Value2 = $False
$ws.Range('H29').Value2 = 86

$ws.Range('A30').Value2 = 'difficult'
$ws.Range('B30').Value2 = 0.2584269662921349
$ws.Range('C30').Value2 = 23
$ws.Range('D30').Value2 = 23
$ws.Range('E30').Value2 = 0
$ws.Range('F30').Value2 = 1
$ws.Range('G30').Value2 = $False
$ws.Range('H30').Value2 = 66

$ws.Range('A31').Value2 = 'size'
$ws.Range('B31').Value2 = 0.2371134020618557
$ws.Range('C31').Value2 = 46
$ws.Range('D31').Value2 = 46
$ws.Range('E31').Value2 = 0
$ws.Range('F31').Value2 = 1
$ws.Range('G31').Value2 = $False
$ws.Range('H31').Value2 = 148

$ws.Range('A32').Value2 = 'item'
$ws.Range('B32').Value2 = 0.2210144927536232
$ws.Range('C32').Value2 = 61
$ws.Range('D32').Value2 = 61
$ws.Range('E32').Value2 = 0
$ws.Range('F32').Value2 = 1
$ws.Range('G32').Value2 = $False
$ws.Range('H32').Value2 = 215

$ws.Range('A33').Value2 = 'hard'
$ws.Range('B33').Value2 = 0.215
$ws.Range('C33').Value2 = 43
$ws.Range('D33').Value2 = 43
$ws.Range('E33').Value2 = 0
$ws.Range('F33').Value2 = 1
$ws.Range('G33').Value2 = $False
$ws.Range('H33').Value2 = 157

$ws.Range('A34').Value2 = 'money'
$ws.Range('B34').Value2 = 0.1962025316455696
$ws.Range('C34').Value2 = 62
$ws.Range('D34').Value2 = 62
$ws.Range('E34').Value2 = 0
$ws.Range('F34').Value2 = 1
$ws.Range('G34').Value2 = $False
$ws.Range('H34').Value2 = 254

$ws.Range('A35').Value2 = 'would'
$ws.Range('B35').Value2 = 0.185459940652819
$ws.Range('C35').Value2 = 125
$ws.Range('D35').Value2 = 125
$ws.Range('E35').Value2 = 0
$ws.Range('F35').Value2 = 1
$ws.Range('G35').Value2 = $False
$ws.Range('H35').Value2 = 549

$ws.Range('A36').Value2 = 'price'
$ws.Range('B36').Value2 = 0.170028818443804
$ws.Range('C36').Value2 = 59
$ws.Range('D36').Value2 = 60
$ws.Range('E36').Value2 = 0.02
$ws.Range('F36').Value2 = 0.98
$ws.Range('G36').Value2 = $True
$ws.Range('H36').Value2 = 288

$ws.Range('A37').Value2 = 'better'
$ws.Range('B37').Value2 = 0.1635514018691589
$ws.Range('C37').Value2 = 35
$ws.Range('D37').Value2 = 35
$ws.Range('E37').Value2 = 0
$ws.Range('F37').Value2 = 1
$ws.Range('G37').Value2 = $False
$ws.Range('H37').Value2 = 179

$ws.Range('A38').Value2 = 'work'
$ws.Range('B38').Value2 = 0.1582278481012658
$ws.Range('C38').Value2 = 50
$ws.Range('D38').Value2 = 50
$ws.Range('E38').Value2 = 0
$ws.Range('F38').Value2 = 1
$ws.Range('G38').Value2 = $False
$ws.Range('H38').Value2 = 266

$ws.Range('A39').Value2 = 'product'
$ws.Range('B39').Value2 = 0.1545253863134658
$ws.Range('C39').Value2 = 70
$ws.Range('D39').Value2 = 71
$ws.Range('E39').Value2 = 0.01
$ws.Range('F39').Value2 = 0.99
$ws.Range('G39').Value2 = $True
$ws.Range('H39').Value2 = 383

$ws.Range('A40').Value2 = 'could'
$ws.Range('B40').Value2 = 0.1401273885350318
$ws.Range('C40').Value2 = 22
$ws.Range('D40').Value2 = 22
$ws.Range('E40').Value2 = 0
$ws.Range('F40').Value2 = 1
$ws.Range('G40').Value2 = $False
$ws.Range('H40').Value2 = 135

$ws.Range('A41').Value2 = 'used'
$ws.Range('B41').Value2 = 0.1314285714285714
$ws.Range('C41').Value2 = 23
$ws.Range('D41').Value2 = 23
$ws.Range('E41').Value2 = 0
$ws.Range('F41').Value2 = 1
$ws.Range('G41').Value2 = $False
$ws.Range('H41').Value2 = 152

$ws.Range('A42').Value2 = '3'
$ws.Range('B42').Value2 = 0.1169354838709677
$ws.Range('C42').Value2 = 29
$ws.Range('D42').Value2 = 29
$ws.Range('E42').Value2 = 0
$ws.Range('F42').Value2 = 1
$ws.Range('G42').Value2 = $False
$ws.Range('H42').Value2 = 219

$ws.Range('A43').Value2 = '2'
$ws.Range('B43').Value2 = 0.1132075471698113
$ws.Range('C43').Value2 = 30
$ws.Range('D43').Value2 = 32
$ws.Range('E43').Value2 = 0.06
$ws.Range('F43').Value2 = 0.9399999999999999
$ws.Range('G43').Value2 = $True
$ws.Range('H43').Value2 = 235

$ws.Range('A44').Value2 = 'use'
$ws.Range('B44').Value2 = 0.09315068493150686
$ws.Range('C44').Value2 = 34
$ws.Range('D44').Value2 = 34
$ws.Range('E44').Value2 = 0
$ws.Range('F44').Value2 = 1
$ws.Range('G44').Value2 = $False
$ws.Range('H44').Value2 = 331

$ws.Range('A45').Value2 = 'little'
$ws.Range('B45').Value2 = 0.08482142857142858
$ws.Range('C45').Value2 = 38
$ws.Range('D45').Value2 = 39
$ws.Range('E45').Value2 = 0.03
$ws.Range('F45').Value2 = 0.97
$ws.Range('G45').Value2 = $True
$ws.Range('H45').Value2 = 410

$ws.Range('A46').Value2 = 'made'
$ws.Range('B46').Value2 = 0.08280254777070063
$ws.Range('C46').Value2 = 26
$ws.Range('D46').Value2 = 29
$ws.Range('E46').Value2 = 0.1
$ws.Range('F46').Value2 = 0.9
$ws.Range('G46').Value2 = $True
$ws.Range('H46').Value2 = 288

$ws.Range('A47').Value2 = 'like'
$ws.Range('B47').Value2 = 0.07425742574257425
$ws.Range('C47').Value2 = 45
$ws.Range('D47').Value2 = 47
$ws.Range('E47').Value2 = 0.04
$ws.Range('F47').Value2 = 0.96
$ws.Range('G47').Value2 = $True
$ws.Range('H47').Value2 = 561

$ws.Range('A48').Value2 = 'buy'
$ws.Range('B48').Value2 = 0.07323943661971831
$ws.Range('C48').Value2 = 26
$ws.Range('D48').Value2 = 26
$ws.Range('E48').Value2 = 0
$ws.Range('F48').Value2 = 1
$ws.Range('G48').Value2 = $False
$ws.Range('H48').Value2 = 329

$ws.Range('A49').Value2 = 'much'
$ws.Range('B49').Value2 = 0.05841121495327103
$ws.Range('C49').Value2 = 25
$ws.Range('D49').Value2 = 31
$ws.Range('E49').Value2 = 0.19
$ws.Range('F49').Value2 = 0.8100000000000001
$ws.Range('G49').Value2 = $True
$ws.Range('H49').Value2 = 403

$ws.Range('A50').Value2 = 'one'
$ws.Range('B50').Value2 = 0.05216284987277354
$ws.Range('C50').Value2 = 41
$ws.Range('D50').Value2 = 49
$ws.Range('E50').Value2 = 0.16
$ws.Range('F50').Value2 = 0.84
$ws.Range('G50').Value2 = $True
$ws.Range('H50').Value2 = 745

$ws.Range('A51').Value2 = 'toy'
$ws.Range('B51').Value2 = 0.04287901990811639
$ws.Range('C51').Value2 = 28
$ws.Range('D51').Value2 = 30
$ws.Range('E51').Value2 = 0.07000000000000001
$ws.Range('F51').Value2 = 0.9299999999999999
$ws.Range('G51').Value2 = $True
$ws.Range('H51').Value2 = 625

$ws.Range('J28:Q28').ClearContents()
